$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Simulador IR 2025")
$ws.Unprotect()
$ws.Range("A1").Value = "test"
Write-Output "done"
